$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.476.85"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").Value = "2.967.53"
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "2.969.78"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.09"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000216"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.99"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "3.440.75"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "60.472.55"
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("D18").Value = "2.965.97"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.95"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.74"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.72"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.60"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.55"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.27%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "54.72"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.31%  "
$ws.Range("E34").Value = "  -6.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "447.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").Value = "3.156.44"
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0768"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.96"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.52%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.97"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("E49").Value = "  -5.29%  "
$ws.Range("D50").Value = "0.0₃0489"
$ws.Range("E50").Value = "  -10.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.61%  "
